$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-72 all get updated from serial date 45172 to 45175
$newDate = [DateTime]::FromOADate(45175)

for ($r = 2; $r -le 72; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
